$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.034.59"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "2.398.17"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'505.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'132.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "2.411.87"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "'0.0969"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("D13").Value = "'4.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").Value = "2.827.40"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "56.960.13"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "'21.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").Value = "2.386.52"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'10.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'309.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "'0.376"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").Value = "'173.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'5.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'17.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "'3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("D40").Value = "'36.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Value = "'0.817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").Value = "'133.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.07%  "
$ws.Range("D44").Value = "'5.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.29%  "
$ws.Range("D45").Value = "'3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "'252.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "'0.567"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'17.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.83%  "
